$wb = $excel.ActiveWorkbook

# --- NewLoanInput sheet ---
$ws1 = $wb.Worksheets.Item("NewLoanInput")

# "Chaithanya 123" -> "chaithanyatest" (shared string reused at B2)
$ws1.Range("B2").Value = "chaithanyatest"

# Insert a new row 7 ("Firstrepaymenton" / date) - pushes old rows 7-17 down to 8-18
$ws1.Rows.Item(7).EntireRow.Insert()
$ws1.Range("A7").Value = "Firstrepaymenton"
$ws1.Range("B7").Value = 42036

# --- Summary sheet: selection only ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A3").Select() | Out-Null

# --- Repayment Schedule sheet: selection only ---
$ws3 = $wb.Worksheets.Item("Repayment Schedule")
$ws3.Range("C8").Select() | Out-Null

# --- NewLoanInput sheet: final selection ---
$ws1.Range("B11").Select() | Out-Null

# --- Transactions sheet: value change + selection (must be last so it stays the active tab) ---
$ws4 = $wb.Worksheets.Item("Transactions")
$ws4.Range("A2").Value = 203
$ws4.Range("H2").Select() | Out-Null
